$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "IdIngeniero" column (column P), shifting the remaining
# columns (mailIng, F_SolicitudServicio, ...) one position to the left.
$ws.Columns("P").Delete()

# Remove the data row (row 2) so only the header row remains.
$ws.Rows(2).Delete()
